# Add a new student row ("MAJAJD" / "Ines") to the assignments table,
# inserted just above the "CHAUVIN" row (the former row 101), pushing the
# rest of the table down by one row and growing the table/autofilter range
# from A1:E115 to A1:E116.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects("Tableau2")

# Insert a blank worksheet row above row 101 (this shifts rows 101-115 down
# to 102-116 and copies the row format, same as Excel's own "Insert" does).
$ws.Rows("101:101").Insert()

# The table's own bookkeeping (ref / autoFilter) doesn't auto-grow from a
# plain row insert, so resize it explicitly to include the new row.
[void]$tbl.Resize($ws.Range("A1:E116"))

# Fill in the new row's values.
$ws.Range("A101").Value = "MAJAJD"
$ws.Range("B101").Value = "Ines"
$ws.Range("C101").Value = "ALL"
$ws.Range("D101").Value = "ALL"
$ws.Range("E101").Value = "ALL"

# Match the saved cursor position recorded in the workbook after the edit.
[void]$ws.Range("H111").Select()
